$d = $word.ActiveDocument

# Locate the paragraph that ends "...MacGrall-Hill" (the last bibliography
# entry we keep) and the paragraph that ends with "...Creative Commons
# Attribution" (the last paragraph we want to remove), then delete
# everything in between (the trailing blank paragraph, the "Ver no
# Jupiter..." paragraph, and the "© 2020 ..." paragraph), including their
# paragraph marks. The paragraph mark that ends "MacGrall-Hill" itself, and
# everything after the "Creative Commons Attribution" paragraph mark, are
# left untouched.

$keepEnd = $d.Content.Duplicate
$keepEnd.Find.Execute("MacGrall-Hill", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$keepEnd.Expand(4) | Out-Null

$removeEnd = $d.Content.Duplicate
$removeEnd.Find.Execute("Creative Commons Attribution", $false, $false, $false, `
                         $false, $false, $true, 1, $false, "", 0) | Out-Null
$removeEnd.Expand(4) | Out-Null

$victim = $d.Range($keepEnd.End, $removeEnd.End)
$victim.Delete()
